$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.2518466851032267
$ws.Range("J2").Value = 0.2518466851032267
$ws.Range("M2").Value = 0.8170803333333333
$ws.Range("N2").Value = 2.451241
$ws.Range("Q2").Value = 0.06982114864399999
$ws.Range("R2").Value = 0.628390337796
$ws.Range("S2").Value = 0.2518466851032267
$ws.Range("T2").Value = 0.2518466851032267

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2538496666666667
$ws.Range("H3").Value = 0.761549
$ws.Range("I3").Value = 0.7481533148967734
$ws.Range("J3").Value = 0.7481533148967733
$ws.Range("M3").Value = 0.8170803333333333
$ws.Range("N3").Value = 2.451241
$ws.Range("Q3").Value = 0.2074155702565556
$ws.Range("R3").Value = 1.866740132309
$ws.Range("S3").Value = 0.7481533148967734
$ws.Range("T3").Value = 0.7481533148967733
